# Add two new exposure-site rows to the table.
# The source data gained two new rows: one becomes the new row 2
# ("400 Dalton Rd, Epping VIC 3076") and one becomes the new row 6
# ("2/44 Hampstead Rd, Maidstone VIC 3012"), with all existing rows
# shifting down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the first new row at position 2 (right after the header row),
# pushing the former row 2 ("195A Stewart St...") and everything below
# it down by one.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()
$ws.Range("A2").Value = "400 Dalton Rd, Epping VIC 3076"
$ws.Range("B2").Value = -37.655759
$ws.Range("C2").Value = 145.032038
$ws.Range("D2").Value = "Whittlesea (C)"

# Insert the second new row at position 6, pushing the former row 6
# ("12-18 Distribution Dr, Truganina...", now at row 7 after the first
# insert) and everything below it down by one more.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).ClearFormats()
$ws.Range("A6").Value = "2/44 Hampstead Rd, Maidstone VIC 3012"
$ws.Range("B6").Value = -37.778539
$ws.Range("C6").Value = 144.875831
$ws.Range("D6").Value = "Maribyrnong (C)"
